# edit.ps1 - apply "order data written into per-table rows" change
# (fills previously-blank placeholder rows / extends sheets with new order
#  rows that were written - including the buggy "Table null"/"Table l"/
#  "Table 0"/"Table 666"/"Table " table-name values described in the commit)
$wb = $excel.ActiveWorkbook

# --- Oto Bento ---
$ws = $wb.Worksheets.Item("Oto Bento")
$ws.Cells.Item(29,1).Value = "Table null"
$ws.Cells.Item(30,1).Value = "Table null"
$ws.Cells.Item(30,2).Value = "Chicken Katsu"
$ws.Cells.Item(30,3).Value = 1.0
$ws.Cells.Item(30,4).Value = "masukkan catatan disini"
$ws.Cells.Item(30,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(31,1).Value = "Table null"
$ws.Cells.Item(31,2).Value = "Chicken Katsu"
$ws.Cells.Item(31,3).Value = 1.0
$ws.Cells.Item(31,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(31,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(32,1).Value = "Table 0"
$ws.Cells.Item(32,2).Value = "Nasi Putih"
$ws.Cells.Item(32,3).Value = 1.0
$ws.Cells.Item(32,4).Value = "masukkan catatan disini"
$ws.Cells.Item(32,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(33,1).Value = "Table 0"
$ws.Cells.Item(33,2).Value = "Chicken Fried Rice"
$ws.Cells.Item(33,3).Value = 1.0
$ws.Cells.Item(33,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(33,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(34,1).Value = "Table 0"
$ws.Cells.Item(34,2).Value = "Chicken Baaga"
$ws.Cells.Item(34,3).Value = 1.0
$ws.Cells.Item(34,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(34,5).Value = "Masukkan catatan disini"

# --- Ootoya ---
$ws = $wb.Worksheets.Item("Ootoya")
$ws.Cells.Item(26,1).Value = "Table 5"
$ws.Cells.Item(26,2).Value = "Jako Gohan"
$ws.Cells.Item(26,3).Value = 4.0
$ws.Cells.Item(26,4).Value = "masukkan catatan disini"
$ws.Cells.Item(26,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(27,1).Value = "Table 5"
$ws.Cells.Item(27,2).Value = "Demi Hamburg"
$ws.Cells.Item(27,3).Value = 8.0
$ws.Cells.Item(27,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(27,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(28,1).Value = "Table 2"
$ws.Cells.Item(28,2).Value = "Hijiki Gohan"
$ws.Cells.Item(28,3).Value = 1.0
$ws.Cells.Item(28,4).Value = "masukkan catatan disini"
$ws.Cells.Item(28,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(29,1).Value = "Table null"
$ws.Cells.Item(29,2).Value = "Beef Curry"
$ws.Cells.Item(29,3).Value = 1.0
$ws.Cells.Item(29,4).Value = "masukkan catatan disini"
$ws.Cells.Item(29,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(30,1).Value = "Table "
$ws.Cells.Item(30,2).Value = "Gokoku Gohan"
$ws.Cells.Item(30,3).Value = 1.0
$ws.Cells.Item(30,4).Value = "masukkan catatan disini"
$ws.Cells.Item(30,5).Value = "Masukkan catatan disini"

# --- Arasseo ---
$ws = $wb.Worksheets.Item("Arasseo")
$ws.Cells.Item(26,1).Value = "Table null"
$ws.Cells.Item(26,2).Value = "Kimchi Tofu"
$ws.Cells.Item(26,2).Style = "Normal"
$ws.Cells.Item(26,3).Value = 1.0
$ws.Cells.Item(26,4).Value = "masukkan catatan disini"
$ws.Cells.Item(26,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(27,1).Value = "Table null"
$ws.Cells.Item(27,2).Value = "Wild Wild Wings"
$ws.Cells.Item(27,3).Value = 2.0
$ws.Cells.Item(27,4).Value = "masukkan catatan disini"
$ws.Cells.Item(27,5).Value = "Masukkan catatan disini"

# --- Kyochon ---
$ws = $wb.Worksheets.Item("Kyochon")
$ws.Cells.Item(26,1).Value = "Table null"
$ws.Cells.Item(26,2).Value = "Kimchi Soup"
$ws.Cells.Item(26,2).Style = "Normal"
$ws.Cells.Item(26,3).Value = 5.0
$ws.Cells.Item(26,4).Value = "masukkan catatan disini"
# Kyochon row 26 col E target is an explicit empty string;
# Excel clears a cell when assigned "", so we leave it unset.
$ws.Cells.Item(27,1).Value = "Table null"
$ws.Cells.Item(27,2).Value = "Mineral Water"
$ws.Cells.Item(27,3).Value = 1.0
$ws.Cells.Item(27,4).Value = "masukkan catatan disini"
$ws.Cells.Item(27,5).Value = "Masukkan catatan disini"

# --- Han Gang ---
$ws = $wb.Worksheets.Item("Han Gang")
$ws.Cells.Item(26,1).Value = "Table null"
$ws.Cells.Item(26,2).Value = "Maeun Tteokpoki"
$ws.Cells.Item(26,2).Style = "Normal"
$ws.Cells.Item(26,3).Value = 1.0
$ws.Cells.Item(26,4).Value = "masukkan catatan disini"
$ws.Cells.Item(26,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(27,1).Value = "Table "
$ws.Cells.Item(27,2).Value = "Ojingeo Bokeum"
$ws.Cells.Item(27,3).Value = 2.0
$ws.Cells.Item(27,4).Value = "masukkan catatan disini"
$ws.Cells.Item(27,5).Value = "Masukkan catatan disini"

# --- PHD ---
$ws = $wb.Worksheets.Item("PHD")
$ws.Cells.Item(40,1).Value = "Table l"
$ws.Cells.Item(40,2).Value = "Boneless Chicken Thigh"
$ws.Cells.Item(40,2).Style = "Normal"
$ws.Cells.Item(40,3).Value = 1.0
$ws.Cells.Item(40,4).Value = "masukkan catatan disini"
$ws.Cells.Item(40,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(41,1).Value = "Table l"
$ws.Cells.Item(41,2).Value = "Cheesy Melt Potato"
$ws.Cells.Item(41,3).Value = 1.0
$ws.Cells.Item(41,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(41,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(42,1).Value = "Table l"
$ws.Cells.Item(42,2).Value = "Beef Spaghetti"
$ws.Cells.Item(42,3).Value = 1.0
$ws.Cells.Item(42,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(42,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(43,1).Value = "Table l"
$ws.Cells.Item(43,2).Value = "Boneless Chicken Thigh"
$ws.Cells.Item(43,3).Value = 1.0
$ws.Cells.Item(43,4).Value = "masukkan catatan disini"
$ws.Cells.Item(43,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(44,1).Value = "Table l"
$ws.Cells.Item(44,2).Value = "Cheesy Melt Potato"
$ws.Cells.Item(44,3).Value = 1.0
$ws.Cells.Item(44,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(44,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(45,1).Value = "Table l"
$ws.Cells.Item(45,2).Value = "Beef Spaghetti"
$ws.Cells.Item(45,3).Value = 1.0
$ws.Cells.Item(45,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(45,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(46,1).Value = "Table "
$ws.Cells.Item(46,2).Value = "Beef Spaghetti"
$ws.Cells.Item(46,3).Value = 3.0
$ws.Cells.Item(46,4).Value = "masukkan catatan disini"
$ws.Cells.Item(46,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(47,1).Value = "Table "
$ws.Cells.Item(47,2).Value = "Cheesy Melt Potato"
$ws.Cells.Item(47,3).Value = 1.0
$ws.Cells.Item(47,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(47,5).Value = "Masukkan catatan disini"

# --- Delicio ---
$ws = $wb.Worksheets.Item("Delicio")
$ws.Cells.Item(26,1).Value = "Table null"
$ws.Cells.Item(26,2).Value = "Wyndham Estate Bin 888 Cabernet Merlot"
$ws.Cells.Item(26,2).Style = "Normal"
$ws.Cells.Item(26,3).Value = 1.0
$ws.Cells.Item(26,4).Value = "masukkan catatan disini"
$ws.Cells.Item(26,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(27,1).Value = "Table null"
$ws.Cells.Item(27,2).Value = "Wyndham Estate Bin 888 Cabernet Merlot"
$ws.Cells.Item(27,3).Value = 1.0
$ws.Cells.Item(27,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(27,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(28,1).Value = "Table null"
$ws.Cells.Item(28,2).Value = "Wyndham Estate Bin 888 Cabernet Merlot"
$ws.Cells.Item(28,3).Value = 1.0
$ws.Cells.Item(28,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(28,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(29,1).Value = "Table 2"
$ws.Cells.Item(29,2).Value = "Paket Box isi 12 mini Pastry"
$ws.Cells.Item(29,3).Value = 1.0
$ws.Cells.Item(29,4).Value = "masukkan catatan disini"
$ws.Cells.Item(29,5).Value = "Masukkan catatan disini"

# --- Bodega ---
$ws = $wb.Worksheets.Item("Bodega")
$ws.Cells.Item(18,1).Value = "Table 2"
$ws.Cells.Item(18,2).Value = "Jack Daniel's"
$ws.Cells.Item(18,2).Style = "Normal"
$ws.Cells.Item(18,3).Value = 2.0
$ws.Cells.Item(18,4).Value = "masukkan catatan disini"
$ws.Cells.Item(18,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(19,1).Value = "Table 666"
$ws.Cells.Item(19,2).Value = "Cockburn's Tawny Port"
$ws.Cells.Item(19,3).Value = 2.0
$ws.Cells.Item(19,4).Value = "masukkan catatan disini"
$ws.Cells.Item(19,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(20,1).Value = "Table 666"
$ws.Cells.Item(20,2).Value = "Kahlua"
$ws.Cells.Item(20,3).Value = 1.0
$ws.Cells.Item(20,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(20,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(21,1).Value = "Table 666"
$ws.Cells.Item(21,2).Value = "Jameson's"
$ws.Cells.Item(21,3).Value = 1.0
$ws.Cells.Item(21,4).Value = "masukkan catatan disini"
$ws.Cells.Item(21,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(22,1).Value = "Table 666"
$ws.Cells.Item(22,2).Value = "Glenfiddich 12"
$ws.Cells.Item(22,3).Value = 1.0
$ws.Cells.Item(22,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(22,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(23,1).Value = "Table 3"
$ws.Cells.Item(23,2).Value = "Jack Daniel's"
$ws.Cells.Item(23,3).Value = 1.0
$ws.Cells.Item(23,4).Value = "masukkan catatan disini"
$ws.Cells.Item(23,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(24,1).Value = "Table 5"
$ws.Cells.Item(24,2).Value = "Grey Goose"
$ws.Cells.Item(24,3).Value = 1.0
$ws.Cells.Item(24,4).Value = "masukkan catatan disini"
$ws.Cells.Item(24,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(25,1).Value = "Table 5"
$ws.Cells.Item(25,2).Value = "Jameson's"
$ws.Cells.Item(25,3).Value = 1.0
$ws.Cells.Item(25,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(25,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(26,1).Value = "Table 5"
$ws.Cells.Item(26,2).Value = "Glenfiddich 12"
$ws.Cells.Item(26,3).Value = 1.0
$ws.Cells.Item(26,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(26,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(27,1).Value = "Table 4"
$ws.Cells.Item(27,2).Value = "Glenfiddich 12"
$ws.Cells.Item(27,3).Value = 1.0
$ws.Cells.Item(27,4).Value = "masukkan catatan disini"
$ws.Cells.Item(27,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(28,1).Value = "Table 4"
$ws.Cells.Item(28,2).Value = "Jack Daniel's"
$ws.Cells.Item(28,3).Value = 1.0
$ws.Cells.Item(28,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(28,5).Value = "Masukkan catatan disini"

# --- Aciap ---
$ws = $wb.Worksheets.Item("Aciap")
$ws.Cells.Item(50,1).Value = "Table 3"
$ws.Cells.Item(50,2).Value = "Bihun Siram"
$ws.Cells.Item(50,2).Style = "Normal"
$ws.Cells.Item(50,3).Value = 1.0
$ws.Cells.Item(50,4).Value = "masukkan catatan disini"
$ws.Cells.Item(50,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(51,1).Value = "Table 3"
$ws.Cells.Item(51,2).Value = "Nasi Putih"
$ws.Cells.Item(51,3).Value = 1.0
$ws.Cells.Item(51,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(51,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(52,1).Value = "Table 3"
$ws.Cells.Item(52,2).Value = "Bihun Yam"
$ws.Cells.Item(52,3).Value = 1.0
$ws.Cells.Item(52,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(52,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(53,1).Value = "Table 3"
$ws.Cells.Item(53,2).Value = "Nasi Siram"
$ws.Cells.Item(53,3).Value = 1.0
$ws.Cells.Item(53,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(53,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(54,1).Value = "Table 3"
$ws.Cells.Item(54,2).Value = "Bihun Bun"
$ws.Cells.Item(54,3).Value = 1.0
$ws.Cells.Item(54,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(54,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(55,1).Value = "Table 3"
$ws.Cells.Item(55,2).Value = "Bihun Bun"
$ws.Cells.Item(55,3).Value = 1.0
$ws.Cells.Item(55,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(55,5).Value = "Masukkan catatan disini"

# --- Kacamata ---
$ws = $wb.Worksheets.Item("Kacamata")
$ws.Cells.Item(71,1).Value = "Table null"
$ws.Cells.Item(71,2).Value = "Nasi Hainam Chasiu Siobak"
$ws.Cells.Item(71,2).Style = "Normal"
$ws.Cells.Item(71,3).Value = 1.0
$ws.Cells.Item(71,4).Value = "masukkan catatan disini"
$ws.Cells.Item(71,5).Value = "Masukkan catatan disini"

# --- Mangkok Putih ---
$ws = $wb.Worksheets.Item("Mangkok Putih")
$ws.Cells.Item(104,1).Value = "Table 2"
$ws.Cells.Item(104,2).Value = "Beef Teriyaki Noodles"
$ws.Cells.Item(104,2).Style = "Normal"
$ws.Cells.Item(104,3).Value = 1.0
$ws.Cells.Item(104,4).Value = "masukkan catatan disini"
$ws.Cells.Item(104,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(105,1).Value = "Table 2"
$ws.Cells.Item(105,2).Value = "Beef Teriyaki"
$ws.Cells.Item(105,3).Value = 13.0
$ws.Cells.Item(105,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(105,5).Value = "Masukkan catatan disini"
$ws.Cells.Item(106,1).Value = "Table 2"
$ws.Cells.Item(106,2).Value = "Beef Katsu"
$ws.Cells.Item(106,3).Value = 1.0
$ws.Cells.Item(106,4).Value = "Masukkan catatan disini"
$ws.Cells.Item(106,5).Value = "Masukkan catatan disini"

# --- Ayam Bakar Ganthari ---
$ws = $wb.Worksheets.Item("Ayam Bakar Ganthari")
$ws.Cells.Item(26,1).Value = "Table null"
$ws.Cells.Item(26,2).Value = "Aqua"
$ws.Cells.Item(26,3).Value = 2.0
$ws.Cells.Item(26,4).Value = "masukkan catatan disini"
$ws.Cells.Item(26,5).Value = "Masukkan catatan disini"
